$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10/11/12 summary values ---
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "39/112"

# Apply mtitleStyle (same style already used by A9) to A10:A12 without
# introducing new style entries - copy formats only from A9.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# --- Remove the third Student/Correct Ans block (columns G:H) entirely ---
$ws.Range("G15:H40").Clear()

# --- Remove the Student/Correct Ans data in columns D:E for rows 19-40 ---
$ws.Range("D19:E40").Clear()

# --- Mark the matching "Student Ans" entries in column A with correctStyle,
#     mirroring the value already present in column B for that row. ---
$ws.Range("B10").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Option A"

$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Option B"

$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Option C"

$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Option B"

$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Option C"

$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Option D"

$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = "Option C"

$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Option D"

$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Option C"

$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = "Option D"

# --- D18 becomes incorrectStyle and now shows "Option B" ---
$ws.Range("D18").Value = "Option B"
$ws.Range("C10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
